$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell D1
$ws.Cells.Item(1, 4).Value = "D219 Location"
$ws.Cells.Item(1, 4).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(1, 4).WrapText = $true

# New data cells in column D
$ws.Cells.Item(2, 4).Value = 4.12
$ws.Cells.Item(2, 4).HorizontalAlignment = -4108  # xlCenter

$ws.Cells.Item(3, 4).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(4, 4).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(5, 4).HorizontalAlignment = -4108  # xlCenter

# New cell E2
$ws.Cells.Item(2, 5).Value = 4.18

# Update selection
$ws.Range("E6").Select()
